$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.052.79"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").Value = "2.302.72"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.28%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.80%  "

$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.982"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.32%  "

$ws.Range("D16").Value = "2.651.49"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("D17").Value = "2.307.96"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "42.014.30"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "287.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.07%  "

$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("E24").Value = "  -1.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.60%  "

$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.72%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.88%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.39"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0883"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.06%  "

$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("E36").Value = "  -9.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.95"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0352"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +21.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.24"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("E44").Value = "  -5.11%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.90%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.15"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.36"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.55%  "
